$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 3839
$ws.Range("C3").Value = 3660
$ws.Range("C4").Value = 2798
$ws.Range("C5").Value = 2010
$ws.Range("C6").Value = 1754
$ws.Range("C7").Value = 787
$ws.Range("C8").Value = 584
$ws.Range("C9").Value = 550
$ws.Range("C10").Value = 526
$ws.Range("B11").Value = "Kitchen & Dining"
$ws.Range("C11").Value = 515

$wb.Save()
